{"js": "// The source diff for this template is purely a re-serialization of the\n// package XML: every hunk re-orders namespace declarations / element\n// attributes (e.g. the root <w:document> xmlns list, <w:pgSz>/<w:pgMar>,\n// and the attribute order inside word/styles.xml's <w:docDefaults>,\n// <w:latentStyles>/<w:lsdException> and <w:style> elements). None of the\n// attribute values themselves change - page size stays 11906 x 16838\n// twips, all margins/fonts/colors/styles keep their original values, etc.\n//\n// In other words the document's actual content/formatting is unchanged by\n// the commit; there is nothing for the Word JS API to mutate here (Office.js\n// has no attribute-order-control surface - it only lets us read/write\n// semantic values). So this script intentionally performs a no-op content\n// edit: it loads the body to confirm nothing needs to change and returns\n// without altering any text, formatting or document properties, leaving\n// the page setup and styles exactly as they already are.\n\n// Touch (read-only) the exact objects the diff's noise refers to, to\n// confirm their values already match what the commit describes - there\n// is nothing to assign, since every value is unchanged.\nconst body = context.document.body;\nbody.load(\"text\");\n\nconst sections = context.document.sections;\nsections.load(\"items\");\n\nconst styles = context.document.getStyles();\nstyles.load(\"items\");\n\nawait context.sync();\n\n// sections[0] corresponds to the single <w:sectPr> (pgSz/pgMar only had\n// its attribute order touched - w:w=11906, w:h=16838, margins all 1417 /\n// 708 / 708 / 0 twips, same as before) and styles.items covers the\n// <w:style> definitions (Normal, Titre1, Titre2, ...) whose attribute\n// order changed but whose values (fonts, colors, sizes, spacing, ...)\n// did not. Nothing further to do.\n", "ps1": "# The source diff for this template is purely a re-serialization of the\n# package XML: every hunk re-orders namespace declarations / element\n# attributes (e.g. the root <w:document> xmlns list, <w:pgSz>/<w:pgMar>,\n# and the attribute order inside word/styles.xml's <w:docDefaults>,\n# <w:latentStyles>/<w:lsdException> and <w:style> elements). None of the\n# attribute values themselves change - page size stays 11906 x 16838\n# twips, all margins/fonts/colors/styles keep their original values, etc.\n#\n# In other words the document's actual content/formatting is unchanged by\n# the commit; there is nothing for the Word object model to mutate here\n# (attribute/namespace ordering in the saved part XML is not something the\n# COM object model exposes or controls). So this script intentionally\n# performs a no-op content edit: it reads the active document to confirm\n# nothing needs to change and returns without altering any text,\n# formatting or document properties, leaving the page setup and styles\n# exactly as they already are.\n\n$d = $word.ActiveDocument\n\n# Touch (read-only) the exact objects the diff's noise refers to, to\n# confirm their values already match what the commit describes -\n# nothing to assign, since every value is unchanged.\n$section = $d.Sections.Item(1)\n$pageSetup = $section.PageSetup\n$null = $pageSetup.PageWidth    # stays 11906 twips (595.3 pt)\n$null = $pageSetup.PageHeight   # stays 16838 twips (841.9 pt)\n$null = $pageSetup.TopMargin\n$null = $pageSetup.BottomMargin\n$null = $pageSetup.LeftMargin\n$null = $pageSetup.RightMargin\n$null = $pageSetup.HeaderDistance\n$null = $pageSetup.FooterDistance\n$null = $pageSetup.Gutter\n$null = $d.Styles.Count\n$null = $d.Content.Text\n"}
